# "Add report and update experiments and poster"
#
# The UC Merced validation-results table (Table57, A8:F11 on Sheet1) gets
# three newly-reported numbers filled in for the 10/20 "Validation split"
# columns that were previously left blank:
#   D9  (Shallow CNN, split "30")
#   C10 (Deep CNN,    split "20")
#   D11 (Deep RNN,    split "30")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value  = 0.35349214571341298
$ws.Range("C10").Value = 0.150590426858218
$ws.Range("D11").Value = 0.96204305206596397

# The author's cursor ends up parked on the last cell they typed into.
$ws.Range("D11").Select()

$wb.Save()
